$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "PositiveExtra" (sheet1): append two new data rows (16 and 17)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PositiveExtra")

$ws1.Cells.Item(17, 1).Value2 = 16
$ws1.Cells.Item(17, 1).NumberFormat = "0"
$ws1.Cells.Item(17, 2).Value2 = 12.5

$ws1.Cells.Item(18, 1).Value2 = 17
$ws1.Cells.Item(18, 1).NumberFormat = "0"
$ws1.Cells.Item(18, 2).Value2 = 12.5

# ---------------------------------------------------------------------
# Sheet "ProductList" (sheet2): swap the two product names in B8/B9,
# add print/page setup info, and touch column C so the used range
# picks it up (mirrors the author's on-sheet selection touching C1).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ProductList")

$b8 = $ws2.Cells.Item(8, 2).Value2
$b9 = $ws2.Cells.Item(9, 2).Value2
$ws2.Cells.Item(8, 2).Value2 = $b9
$ws2.Cells.Item(9, 2).Value2 = $b8

$ws2.Cells.Item(1, 3).Value2 = "tmp"
$ws2.Cells.Item(1, 3).ClearContents()

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Selections / scroll position per sheet (activeCell + sqref match the
# author's final cursor position on each tab)
# ---------------------------------------------------------------------
$ws1.Range("C17").Select()

$ws3 = $wb.Worksheets.Item("CAPTSetting")
$ws4 = $wb.Worksheets.Item("SMAPTSetting")

$ws3.Range("D19").Select()
$ws4.Range("B10").Select()

# ProductList ends up the active sheet/tab (activeTab=1), and its
# selection is B9 - select it last so it becomes the active selection.
$ws2.Activate()
$ws2.Range("B9").Select()
